# Natmi following Dr Hou advice
# Expand the LR-pair result table from 6 data rows (rows 2-7) to 12 data
# rows (rows 2-13), recomputing all of the specificity / weight columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster, Target cluster, Ligand-expressing cells, Ligand avg expr,
# Ligand total expr, Ligand derived specificity (avg), Ligand derived
# specificity (total), Receptor-expressing cells, Receptor detection rate,
# Receptor avg expr, Receptor total expr, Receptor derived specificity (avg),
# Receptor derived specificity (total), Edge avg expr weight, Edge total
# expr weight, Edge avg expr derived specificity, Edge total expr derived
# specificity
$rows = @(
    @("ECs",    "FAPs", 2, 195.0792385,        390.158477,  0.2640605522989327,  0.1982306263353075, 2, 0.6666666666666666, 0.1074876666666667, 0.322463,  0.111695032377957,  0.1586806295398324, 20.96861216147517,  125.811672968851,  0.02949425193877049, 0.03145536058096188),
    @("ECs",    "sCs",  2, 195.0792385,        390.158477,  0.2640605522989327,  0.1982306263353075, 2, 1,                   0.8548439999999999,  1.709688,  0.888304967622043,  0.8413193704601676, 166.762316556294,   667.049266225176,  0.2345663003601622,  0.1667752657543456),
    @("FAPs",   "FAPs", 3, 62.40792233333334,  187.223767,  0.08447577797556809, 0.09512412720758515, 2, 0.6666666666666666, 0.1074876666666667, 0.322463,  0.111695032377957,  0.1586806295398324, 6.708081953124556,  60.372737578121,   0.009435524756134185, 0.01509435638972671),
    @("FAPs",   "sCs",  3, 62.40792233333334,  187.223767,  0.08447577797556809, 0.09512412720758515, 2, 1,                   0.8548439999999999,  1.709688,  0.888304967622043,  0.8413193704601676, 53.349037959116,    320.094227754696,  0.0750402532194339,  0.08002977081785843),
    @("M1",     "FAPs", 3, 142.8621113333333,  428.586334,  0.1933791023142199,  0.2177549443006804, 2, 0.6666666666666666, 0.1074876666666667, 0.322463,  0.111695032377957,  0.1586806295398324, 15.35591500229356,  138.203235020642,  0.02159948509420705, 0.03455349164704312),
    @("M1",     "sCs",  3, 142.8621113333333,  428.586334,  0.1933791023142199,  0.2177549443006804, 2, 1,                   0.8548439999999999,  1.709688,  0.888304967622043,  0.8413193704601676, 122.124818700632,   732.7489122037919, 0.1717796172200129,  0.1832014526536373),
    @("M2",     "FAPs", 3, 141.6168416666667,  424.850525,  0.1916934970264942,  0.2158568649262854, 2, 0.6666666666666666, 0.1074876666666667, 0.322463,  0.111695032377957,  0.1586806295398324, 15.22206387145278,  136.998574843075,  0.02141121135701807, 0.03425230321699754),
    @("M2",     "sCs",  3, 141.6168416666667,  424.850525,  0.1916934970264942,  0.2158568649262854, 2, 1,                   0.8548439999999999,  1.709688,  0.888304967622043,  0.8413193704601676, 121.0603073977,     726.3618443862,    0.1702822856694761,  0.1816045617092878),
    @("Neutro", "FAPs", 3, 143.783834,         431.351502,  0.1946267522348261,  0.2191598631141254, 2, 0.6666666666666666, 0.1074876666666667, 0.322463,  0.111695032377957,  0.1586806295398324, 15.45498882104733,  139.094899389426,  0.02173884139248552, 0.03477642504881293),
    @("Neutro", "sCs",  3, 143.783834,         431.351502,  0.1946267522348261,  0.2191598631141254, 2, 1,                   0.8548439999999999,  1.709688,  0.888304967622043,  0.8413193704601676, 122.912747791896,   737.4764867513759, 0.1728879108423406,  0.1843834380653125),
    @("sCs",    "FAPs", 2, 53.01711450000001,  106.034229,  0.07176431814995911, 0.05387357411601602, 2, 0.6666666666666666, 0.1074876666666667, 0.322463,  0.111695032377957,  0.1586806295398324, 5.6986859310045,    34.192115586027,   0.00801571783934169, 0.008548692656290244),
    @("sCs",    "sCs",  2, 53.01711450000001,  106.034229,  0.07176431814995911, 0.05387357411601602, 2, 1,                   0.8548439999999999,  1.709688,  0.888304967622043,  0.8413193704601676, 45.321362227638,    181.285448910552,  0.06374860031061742, 0.04532488145972578)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value  = $row[0]   # A: Sending cluster
    $ws.Cells.Item($r, 2).Value  = "Gnai2"   # B: Ligand symbol
    $ws.Cells.Item($r, 3).Value  = "Lpar3"   # C: Receptor symbol
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D: Target cluster
    $ws.Cells.Item($r, 5).Value  = $row[2]   # E: Ligand-expressing cells
    $ws.Cells.Item($r, 6).Value  = 1         # F: Ligand detection rate
    $ws.Cells.Item($r, 7).Value  = $row[3]   # G: Ligand average expression value
    $ws.Cells.Item($r, 8).Value  = $row[4]   # H: Ligand total expression value
    $ws.Cells.Item($r, 9).Value  = $row[5]   # I: Ligand derived specificity (avg)
    $ws.Cells.Item($r, 10).Value = $row[6]   # J: Ligand derived specificity (total)
    $ws.Cells.Item($r, 11).Value = $row[7]   # K: Receptor-expressing cells
    $ws.Cells.Item($r, 12).Value = $row[8]   # L: Receptor detection rate
    $ws.Cells.Item($r, 13).Value = $row[9]   # M: Receptor average expression value
    $ws.Cells.Item($r, 14).Value = $row[10]  # N: Receptor total expression value
    $ws.Cells.Item($r, 15).Value = $row[11]  # O: Receptor derived specificity (avg)
    $ws.Cells.Item($r, 16).Value = $row[12]  # P: Receptor derived specificity (total)
    $ws.Cells.Item($r, 17).Value = $row[13]  # Q: Edge average expression weight
    $ws.Cells.Item($r, 18).Value = $row[14]  # R: Edge total expression weight
    $ws.Cells.Item($r, 19).Value = $row[15]  # S: Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value = $row[16]  # T: Edge total expression derived specificity
    $r = $r + 1
}
